# Selenium.maven.demo/TestData/TestDataSheet.xlsx - "Added New project and codes"
#
# 1. Change the selection on the "Sheet2" tab from F12 to C1.
# 2. Add a brand-new "Sheet1" tab (with sheetId 4) after "Sheet3", containing
#    an email address + password pair, each hyperlinked, styled with the
#    workbook's existing "Hyperlink" cell style, and make it the active tab
#    (selection C13).

$wb = $excel.ActiveWorkbook

# --- 1. Sheet2 tab: move the selection from F12 to C1 -----------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("C1").Select()

# --- 2. Add the new "Sheet1" tab after the last existing sheet --------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Sheet1"

$newSheet.Range("A1").Value = "pranjal.handique@telusinternational.com"
$newSheet.Range("B1").Value = "Welcome@123"

$newSheet.Hyperlinks.Add($newSheet.Range("A1"), "mailto:pranjal.handique@telusinternational.com")
$newSheet.Hyperlinks.Add($newSheet.Range("B1"), "mailto:pranjal.handique@telusinternational.com")

# Match the look of the existing hyperlink cells on Sheet3 / match column widths
$newSheet.Columns.Item(1).ColumnWidth = 38.43
$newSheet.Columns.Item(2).ColumnWidth = 13.71
$newSheet.Range("A1:B1").Style = "Hyperlink"

# Make "Sheet1" the active tab with C13 selected (this also clears
# tabSelected on whichever sheet was previously active).
$newSheet.Range("C13").Select()
